$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetB")

# Update the "Notes" (F2) and "Task" (D2) cells for row 2 with the new
# comment-derived text. F2 is set first so its shared string is registered
# ahead of D2's, matching the authoring order recorded in the workbook.
$ws.Range("F2").Value = "AsanaAPISync_v0.0.2alpha: Dynamic fetching test 2"
$ws.Range("D2").Value = "Meeting: Intake"

# Widen column D to comfortably fit the new, longer task text.
$ws.Range("D1").ColumnWidth = 33

# Leave the selection on the cell that was just edited.
$ws.Range("D2").Select()
